$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$desc = "%98 pamuk içeriği ile nefes alabilen yapıda, cildinize nazik dokunuşlar sunar.Bağcıklı kapama şekliyle kişisel zevke göre ayarlama imkanı taşır.Lastikli bel detayı sayesinde  rahatlığından ödün vermez.34-46 Beden aralığı mevcuttur."

# Row 8 currently holds the single "Bel Paça Lastik  Pantolon" product.
# It becomes the "Antrasit" variant, and two new rows (Gri, Füme) are
# appended below it, each reusing price / category / description.
# Product names are entered first (A8, A9, A10), then the rest of the
# columns are filled in, matching how the strings were typed originally.

$ws.Cells.Item(8, 1).Value = "Bel Paça Lastik  Pantolon Antrasit"
$ws.Cells.Item(9, 1).Value = "Bel Paça Lastik  Pantolon Gri"
$ws.Cells.Item(10, 1).Value = "Bel Paça Lastik  Pantolon Füme"

$ws.Cells.Item(8, 2).Value = "350 Tl"
$ws.Cells.Item(8, 3).Value = "Jeans"
$ws.Cells.Item(8, 4).Value = "ANTRASİT.jpg"
$ws.Cells.Item(8, 5).Value = $desc

$ws.Cells.Item(9, 2).Value = "350 Tl"
$ws.Cells.Item(9, 3).Value = "Jeans"
$ws.Cells.Item(9, 4).Value = "GRİ.jpg"
$ws.Cells.Item(9, 5).Value = $desc

$ws.Cells.Item(10, 2).Value = "350 Tl"
$ws.Cells.Item(10, 3).Value = "Jeans"
$ws.Cells.Item(10, 4).Value = "MAVİİ.jpg"
$ws.Cells.Item(10, 5).Value = $desc

$ws.Cells.Item(11, 4).Select()
